# AFDP-157 - Add Access Control List to SOLR documents - apply assignment
# and data access control rules to case files.
#
# Adds two new rows to the "Assignment Rules" rule table on Sheet1 mirroring
# the existing COMPLAINT default-assignee / default-access rules, but for
# CASE_FILE objects:
#   Row 20: Case File - Default assignee  (assignee, ann-acm)
#   Row 21: Case File - Default access    (*, *)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 20 - Case File default assignee rule
$ws.Range("B20").Value = "Case File – Default assignee"
$ws.Range("C20").Value = "CASE_FILE"
$ws.Range("D20").Value = "participants.?[participantType == 'assignee'].isEmpty()"
$ws.Range("G20").Value = "assignee, ann-acm"

# Row 21 - Case File default access rule
$ws.Range("B21").Value = "Case File – Default access"
$ws.Range("C21").Value = "CASE_FILE"
$ws.Range("D21").Value = "participants.?[participantType == '*'].isEmpty()"
$ws.Range("G21").Value = "*, *"

# The D column cells on these two new rule rows pick up the TRUE/FALSE
# boolean number format, matching the other condition cells in this table.
$ws.Range("D20:D21").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

# Slightly taller rows to fit the new content, matching the rest of the edit.
$ws.Range("A19:A21").RowHeight = 13.8

# Move the active selection to reflect where the edit was made.
$ws.Range("G22").Select() | Out-Null
